$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The update adds two new weekly price rows for "Melón" (Tuna) at what is
# currently row 70, pushing the existing rows 70-124 down to 72-126.
# Insert two whole rows at row 70 (shifts 70:124 -> 72:126, carrying
# formatting - e.g. the date style on column D - down with them).
$ws.Rows("70:71").Insert()

# Populate the newly inserted row 70.
$ws.Range("A70").Value = 8
$ws.Range("B70").Value = "Terminal La Palmera de La Serena"
$ws.Range("C70").Value = "Coquimbo"
$ws.Range("D70").Value = 44978
$ws.Range("E70").Value = 4
$ws.Range("F70").Value = 100112027
$ws.Range("G70").Value = "Melón"
$ws.Range("H70").Value = "Tuna"
$ws.Range("I70").Value = "Extra"
$ws.Range("J70").Value = 1800
$ws.Range("K70").Value = 1400
$ws.Range("L70").Value = 1500
$ws.Range("M70").Value = 1450
$ws.Range("N70").Value = "`$/unidad"
$ws.Range("O70").Value = "Región de O'Higgins"
$ws.Range("P70").Value = 1450
$ws.Range("Q70").Value = 1
$ws.Range("R70").Value = "Hortaliza"

# Populate the newly inserted row 71.
$ws.Range("A71").Value = 8
$ws.Range("B71").Value = "Terminal La Palmera de La Serena"
$ws.Range("C71").Value = "Coquimbo"
$ws.Range("D71").Value = 44978
$ws.Range("E71").Value = 4
$ws.Range("F71").Value = 100112027
$ws.Range("G71").Value = "Melón"
$ws.Range("H71").Value = "Tuna"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 1300
$ws.Range("K71").Value = 1100
$ws.Range("L71").Value = 1200
$ws.Range("M71").Value = 1150
$ws.Range("N71").Value = "`$/unidad"
$ws.Range("O71").Value = "Región de O'Higgins"
$ws.Range("P71").Value = 1150
$ws.Range("Q71").Value = 1
$ws.Range("R71").Value = "Hortaliza"
